# Updates to Rank and other scripts
#
# 1. Remove the two "Lake Creek Methow" rows (01 and 02).
# 2. Remove the two "Wolf Creek" rows (01 and 02).
# 3. Insert a new reach "Entiat River Lake 04" as the new row 2 (just above
#    "Entiat River Potato 07"), with its own scored data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
$lastCol = $colLetters[$colLetters.Length - 1]

function Get-LastRow {
    return $ws.Range("A1048576").End(-4162).Row   # xlUp = -4162
}

function Delete-ReachRow([string]$reachName) {
    $lastRow = Get-LastRow
    $searchRange = $ws.Range("A2:A" + $lastRow)
    $found = $searchRange.Find($reachName)
    if ($found -ne $null) {
        $ws.Rows($found.Row).Delete()
    }
}

# --- 1 & 2: delete the four obsolete rows ---
Delete-ReachRow "Lake Creek Methow 01"
Delete-ReachRow "Lake Creek Methow 02"
Delete-ReachRow "Wolf Creek 01"
Delete-ReachRow "Wolf Creek 02"

# --- 3: insert the new "Entiat River Lake 04" row right above "Entiat River Potato 07" ---
$lastRow = Get-LastRow
$anchor = $ws.Range("A2:A" + $lastRow).Find("Entiat River Potato 07")
$newRowNum = $anchor.Row

$ws.Rows($newRowNum).Insert()
$ws.Rows($newRowNum).ClearFormats()

$values = @(
    "Entiat River Lake 04",
    "Entiat",
    "Entiat River-Lake Creek",
    "yes",
    "yes",
    "yes",
    5,
    5,
    5,
    5,
    1,
    3,
    5,
    3,
    1,
    5,
    5,
    5,
    5,
    33,
    0.7333333333333333,
    5,
    3,
    "Cover-Wood,PoolQuantity&Quality",
    "Flow-SummerBaseFlow,Off-Channel-Side-Channels",
    "Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Side-Channels,PoolQuantity&Quality"
)

for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $cellRef = $colLetters[$i] + $newRowNum
    $ws.Range($cellRef).Value = $values[$i]
}

$usedLastRow = Get-LastRow
Write-Output ("Final last row: " + $usedLastRow)
